$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.346.19'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '3.688.63'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '681.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.65%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('E12').Value = '  -2.83%  '
$ws.Range('D13').Value = '4.310.49'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').Value = '3.699.56'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '69.332.33'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('E17').Value = '  +1.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.656'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '80.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '3.835.08'
$ws.Range('E26').Value = '  -4.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.53%  '
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('E32').Value = '  -2.98%  '
$ws.Range('E33').Value = '  +0.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '3.675.99'
$ws.Range('E35').Value = '  +0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.156'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.75%  '
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  -2.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '170.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.943'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.63'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.69%  '
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E48').Value = '  -2.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000276'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.79'
$ws.Range('D51').Style = 'Normal'
